# Palaverimuistio ja työajan lisäys
# Fill in the first empty row of the time-tracking table with a new entry.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Row 1 is the header; rows 2-8 already contain the 7 existing entries,
# so row 9 is the first blank row to populate.
$row = $table.Rows.Item(9)

$row.Cells.Item(1).Range.Text = "10.02.23"
$row.Cells.Item(2).Range.Text = "0,5"
$row.Cells.Item(3).Range.Text = "Sprintti tapaaminen"
